$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 56 hold a "last changed" date serial that was
# bumped by one day (45184 -> 45185) during the automatic update.
for ($row = 2; $row -le 56; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
